$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (Lecture 4 / link) below the existing table.
$ws.Range("A9").Value = "Лекция 4"
$ws.Range("B9").Value = "https://youtu.be/pwpuLsnDxpg"

# Match the saved selection state from the authored workbook.
$ws.Range("A9").Select()
